# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") on Sheet1 holds the literal string "6-30-2012-13"
# for every data row (rows 2-31). Re-write it as "2013-06-30" while
# keeping the cells as plain text (not auto-converted to a date serial)
# and without leaving any new per-cell formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")

# Force text interpretation so Excel doesn't reinterpret the ISO-style
# string as a date serial number when we assign it below.
$dateRange.NumberFormat = "@"
$dateRange.Value = "2013-06-30"

# Restore the default "Normal" style so the cells end up exactly like
# before (no explicit style index), matching the original formatting.
$dateRange.Style = "Normal"
